# Add a new "Height" property/column for NPC data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header in column AB (28), after the existing last column AA (27).
$ws.Cells.Item(1, 28).Value = "Height"

# Fill the new Height column with a default value of 2 for every data row.
$lastRow = 21
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 28).Value = 2
}

# Match the column width used by the neighbouring column AA.
$ws.Range("AB1:AB21").ColumnWidth = 14.75

# Reflect the selection the author left on the new column in the saved view.
$ws.Range("AB2:AB21").Select()
